$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BSNE")
$ws.Rows.Item(122).Select()
$ws.Rows.Item(122).Insert()
Write-Host "done"
